$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("case_1")
$ws.Range("DD2").Value = 0.0002251
$ws.Range("DE2").Value = 0.001169
$ws.Range("DF2").Value = 0.001483
$ws.Range("DG2").Value = 0.0006431
$ws.Range("DJ2").Value = 0.0002136
$ws.Range("DK2").Value = 0.001246
$ws.Range("DL2").Value = 0.001471
$ws.Range("DM2").Value = 0.0006898
$ws.Range("DP2").Value = 0.0005961
$ws.Range("DQ2").Value = 0.005942
$ws.Range("DR2").Value = 0.009937
$ws.Range("DS2").Value = 0.004479
$ws.Range("DV2").Value = 0.0005432
$ws.Range("DW2").Value = 0.004338
$ws.Range("DX2").Value = 0.006215
$ws.Range("DY2").Value = 0.003043
$ws.Range("EB2").Value = 0.0003913
$ws.Range("EC2").Value = 0.002518
$ws.Range("ED2").Value = 0.003136
$ws.Range("EE2").Value = 0.001484
$ws.Range("DD3").Value = 0.00007014
$ws.Range("DE3").Value = 0.0004082
$ws.Range("DF3").Value = 0.0006963
$ws.Range("DG3").Value = 0.0002639
$ws.Range("DJ3").Value = 0.0002273
$ws.Range("DK3").Value = 0.001183
$ws.Range("DL3").Value = 0.001842
$ws.Range("DM3").Value = 0.0008363
$ws.Range("DP3").Value = 0.000723
$ws.Range("DQ3").Value = 0.004759
$ws.Range("DR3").Value = 0.02939
$ws.Range("DS3").Value = 0.006789
$ws.Range("DV3").Value = 0.0009963999999999999
$ws.Range("DW3").Value = 0.005813
$ws.Range("DX3").Value = 0.05211
$ws.Range("DY3").Value = 0.009823
$ws.Range("EB3").Value = 0.0001699
$ws.Range("EC3").Value = 0.0009661
$ws.Range("ED3").Value = 0.001679
$ws.Range("EE3").Value = 0.0006628

$ws = $wb.Worksheets.Item("case_2")
$ws.Range("DD2").Value = 0.00007873
$ws.Range("DE2").Value = 0.002901
$ws.Range("DF2").Value = 0.003468
$ws.Range("DG2").Value = 0.001687
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = 0.002617
$ws.Range("DL2").Value = 0.003009
$ws.Range("DM2").Value = 0.001477
$ws.Range("DP2").Value = 0.0001019
$ws.Range("DQ2").Value = 0.007505
$ws.Range("DR2").Value = 0.01052
$ws.Range("DS2").Value = 0.006297
$ws.Range("DV2").Value = 0.0002161
$ws.Range("DW2").Value = 0.006427
$ws.Range("DX2").Value = 0.008392
$ws.Range("DY2").Value = 0.004939
$ws.Range("EB2").Value = 0.00004101
$ws.Range("EC2").Value = 0.004813
$ws.Range("ED2").Value = 0.005891
$ws.Range("EE2").Value = 0.00305
$ws.Range("DD3").Value = 0.0001943
$ws.Range("DE3").Value = 0.001957
$ws.Range("DF3").Value = 0.003525
$ws.Range("DG3").Value = 0.001291
$ws.Range("DJ3").Value = 0.0002148
$ws.Range("DK3").Value = 0.003115
$ws.Range("DL3").Value = 0.005161
$ws.Range("DM3").Value = 0.002257
$ws.Range("DP3").Value = 0.000681
$ws.Range("DQ3").Value = 0.007752
$ws.Range("DR3").Value = 0.08286
$ws.Range("DS3").Value = 0.01372
$ws.Range("DV3").Value = 0.0005789
$ws.Range("DW3").Value = 0.02644
$ws.Range("DX3").Value = 0.1395
$ws.Range("DY3").Value = 0.02314
$ws.Range("EB3").Value = 0.0001511
$ws.Range("EC3").Value = 0.002988
$ws.Range("ED3").Value = 0.005618
$ws.Range("EE3").Value = 0.002287

$ws = $wb.Worksheets.Item("case_3")
$ws.Range("E2").Value = 0.04024
$ws.Range("F2").Value = 0.2583
$ws.Range("G2").Value = 0.6158
$ws.Range("H2").Value = 0.8424
$ws.Range("I2").Value = 0.3207
$ws.Range("E3").Value = 0.1027
$ws.Range("F3").Value = 0.5762
$ws.Range("G3").Value = 0.9725
$ws.Range("I3").Value = 0.543

$ws = $wb.Worksheets.Item("case_1_worst_case")
$ws.Range("F2").Value = 0.0002136
$ws.Range("G2").Value = 0.001246
$ws.Range("H2").Value = 0.001471
$ws.Range("I2").Value = 0.0006898
$ws.Range("M2").Value = 0.0005961
$ws.Range("N2").Value = 0.005942
$ws.Range("O2").Value = 0.009937
$ws.Range("P2").Value = 0.004479
$ws.Range("F3").Value = 0.0002273
$ws.Range("G3").Value = 0.001183
$ws.Range("H3").Value = 0.001842
$ws.Range("I3").Value = 0.0008363
$ws.Range("M3").Value = 0.0009963999999999999
$ws.Range("N3").Value = 0.005813
$ws.Range("O3").Value = 0.05211
$ws.Range("P3").Value = 0.009823

$ws = $wb.Worksheets.Item("case_2_worst_case")
$ws.Range("F2").Value = 0.00007873
$ws.Range("G2").Value = 0.002901
$ws.Range("H2").Value = 0.003468
$ws.Range("I2").Value = 0.001687
$ws.Range("M2").Value = 0.0001019
$ws.Range("N2").Value = 0.007505
$ws.Range("O2").Value = 0.01052
$ws.Range("P2").Value = 0.006297
$ws.Range("F3").Value = 0.0002148
$ws.Range("G3").Value = 0.003115
$ws.Range("H3").Value = 0.005161
$ws.Range("I3").Value = 0.002257
$ws.Range("M3").Value = 0.0005789
$ws.Range("N3").Value = 0.02644
$ws.Range("O3").Value = 0.1395
$ws.Range("P3").Value = 0.02314
